# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (use K instead of Strike#, regen std/mean, calc and write s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G2"  = 1
    "G4"  = 3
    "G5"  = 4
    "G6"  = 0
    "G7"  = 1
    "G8"  = 0
    "G9"  = 1
    "G11" = 1
    "G12" = 1
    "G13" = 2
    "G14" = 1
    "G15" = 2
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
